$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $escaped = $text.Replace('"', '""')
    $cell.Formula = '="" & "' + $escaped + '"'
    $excel.Calculate()
    $cell.Copy() | Out-Null
    $cell.PasteSpecial(-4163) | Out-Null
}

$ws.Cells.Item(2, 4).Value = '22.460.68'
$ws.Cells.Item(2, 5).Value = '  -0.10%  '

$ws.Cells.Item(3, 4).Value = '1.570.80'
$ws.Cells.Item(3, 5).Value = '  -0.23%  '

Set-TextValue $ws.Cells.Item(4, 4) '1.002'
$ws.Cells.Item(4, 5).Value = '  +0.06%  '

$ws.Cells.Item(5, 5).Value = '  +0.01%  '

Set-TextValue $ws.Cells.Item(6, 4) '286.22'
$ws.Cells.Item(6, 5).Value = '  -2.01%  '

Set-TextValue $ws.Cells.Item(7, 4) '0.3651'
$ws.Cells.Item(7, 5).Value = '  -1.92%  '

Set-TextValue $ws.Cells.Item(8, 4) '48.08'
$ws.Cells.Item(8, 5).Value = '  -3.70%  '

Set-TextValue $ws.Cells.Item(9, 4) '0.3328'
$ws.Cells.Item(9, 5).Value = '  -2.37%  '

Set-TextValue $ws.Cells.Item(10, 4) '1.127'
$ws.Cells.Item(10, 5).Value = '  -2.09%  '

Set-TextValue $ws.Cells.Item(11, 4) '0.07423'
$ws.Cells.Item(11, 5).Value = '  -1.80%  '

$ws.Cells.Item(12, 5).Value = '  +0.06%  '

Set-TextValue $ws.Cells.Item(13, 4) '20.83'
$ws.Cells.Item(13, 5).Value = '  -2.17%  '

Set-TextValue $ws.Cells.Item(14, 4) '5.971'
$ws.Cells.Item(14, 5).Value = '  -1.36%  '

Set-TextValue $ws.Cells.Item(15, 4) '6.908'
$ws.Cells.Item(15, 5).Value = '  -1.01%  '

$ws.Cells.Item(16, 4).Value = '1.566.98'
$ws.Cells.Item(16, 5).Value = '  -0.32%  '

$ws.Cells.Item(17, 5).Value = '  -1.84%  '

Set-TextValue $ws.Cells.Item(18, 4) '88.02'
$ws.Cells.Item(18, 5).Value = '  -3.34%  '

Set-TextValue $ws.Cells.Item(19, 4) '0.06725'
$ws.Cells.Item(19, 5).Value = '  -0.45%  '

$ws.Cells.Item(20, 5).Value = '  +0.04%  '

Set-TextValue $ws.Cells.Item(21, 4) '6.371'
$ws.Cells.Item(21, 5).Value = '  +1.03%  '

Set-TextValue $ws.Cells.Item(22, 4) '16.40'
$ws.Cells.Item(22, 5).Value = '  +0.11%  '

$ws.Cells.Item(23, 5).Value = '  -1.05%  '

$ws.Cells.Item(24, 4).Value = '22.450.08'
$ws.Cells.Item(24, 5).Value = '  -0.12%  '

Set-TextValue $ws.Cells.Item(25, 4) '2.385'
$ws.Cells.Item(25, 5).Value = '  +0.53%  '

Set-TextValue $ws.Cells.Item(26, 4) '2.614'
$ws.Cells.Item(26, 5).Value = '  -0.67%  '

Set-TextValue $ws.Cells.Item(27, 4) '151.49'
$ws.Cells.Item(27, 5).Value = '  +1.60%  '

Set-TextValue $ws.Cells.Item(28, 4) '19.51'
$ws.Cells.Item(28, 5).Value = '  -2.62%  '

$ws.Cells.Item(29, 5).Value = '  -0.99%  '

Set-TextValue $ws.Cells.Item(30, 4) '124.28'
$ws.Cells.Item(30, 5).Value = '  -1.09%  '

$ws.Cells.Item(31, 4).Value = '1.747.45'
$ws.Cells.Item(31, 5).Value = '  +0.05%  '

Set-TextValue $ws.Cells.Item(32, 4) '1.039'
$ws.Cells.Item(32, 5).Value = '  -3.79%  '

Set-TextValue $ws.Cells.Item(33, 4) '6.132'
$ws.Cells.Item(33, 5).Value = '  -1.74%  '

Set-TextValue $ws.Cells.Item(34, 4) '1.998'
$ws.Cells.Item(34, 5).Value = '  -0.68%  '

Set-TextValue $ws.Cells.Item(35, 4) '9.765'
$ws.Cells.Item(35, 5).Value = '  -0.95%  '

Set-TextValue $ws.Cells.Item(36, 4) '0.08256'
$ws.Cells.Item(36, 5).Value = '  -1.51%  '

Set-TextValue $ws.Cells.Item(37, 4) '0.02420'
$ws.Cells.Item(37, 5).Value = '  -2.85%  '

Set-TextValue $ws.Cells.Item(38, 4) '0.2245'
$ws.Cells.Item(38, 5).Value = '  -2.54%  '

Set-TextValue $ws.Cells.Item(39, 4) '0.06454'
$ws.Cells.Item(39, 5).Value = '  -1.70%  '

Set-TextValue $ws.Cells.Item(40, 4) '5.410'
$ws.Cells.Item(40, 5).Value = '  -1.04%  '

$ws.Cells.Item(41, 5).Value = '  -3.24%  '

Set-TextValue $ws.Cells.Item(42, 4) '11.28'
$ws.Cells.Item(42, 5).Value = '  -0.92%  '

Set-TextValue $ws.Cells.Item(43, 4) '0.6258'
$ws.Cells.Item(43, 5).Value = '  +0.00%  '

$ws.Cells.Item(44, 2).Value = 'EnergySwap'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Cells.Item(44, 4) '13.82'
$ws.Cells.Item(44, 5).Value = '  -1.18%  '

$ws.Cells.Item(45, 2).Value = 'Decentraland'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextValue $ws.Cells.Item(45, 4) '0.6066'
$ws.Cells.Item(45, 5).Value = '  +3.68%  '

$ws.Cells.Item(46, 2).Value = 'PancakeSwap'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue $ws.Cells.Item(46, 4) '3.743'
$ws.Cells.Item(46, 5).Value = '  -1.89%  '

$ws.Cells.Item(47, 2).Value = 'NEARProtocol'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Cells.Item(47, 4) '2.041'
$ws.Cells.Item(47, 5).Value = '  -1.90%  '

$ws.Cells.Item(48, 2).Value = 'Quant'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue $ws.Cells.Item(48, 4) '123.92'
$ws.Cells.Item(48, 5).Value = '  -4.89%  '

$ws.Cells.Item(49, 2).Value = 'EOS'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
Set-TextValue $ws.Cells.Item(49, 4) '1.219'
$ws.Cells.Item(49, 5).Value = '  +0.07%  '

$ws.Cells.Item(50, 2).Value = 'Cronos'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws.Cells.Item(50, 4) '0.07219'
$ws.Cells.Item(50, 5).Value = '  -1.68%  '

$ws.Cells.Item(51, 2).Value = 'Aave'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Cells.Item(51, 4) '76.07'
$ws.Cells.Item(51, 5).Value = '  -0.96%  '
